# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (fund holdings detail) between the
# existing "总计" (summary) and "2021-Q4" sheets, fills it with the
# Q3-2022 fund table, and records the new quarter in the "总计" summary
# sheet (shifting the previous "2021-Q4" summary row down one row).

function Set-TextValue($sheet, $addr, $val) {
    # Forces a numeric-looking string (e.g. "3.95", "012920") to be stored
    # as text rather than being auto-coerced to a number by .Value, while
    # leaving the cell's style untouched (NumberFormat is reset to General
    # / default style right after the write).
    $sheet.Range($addr).NumberFormat = "@"
    $sheet.Range($addr).Value = $val
    $sheet.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$sheetTotal = $wb.Worksheets.Item("总计")
$sheet2021 = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet, placed right before "2021-Q4" so the
#    final tab order is: 总计, 2022-Q3, 2021-Q4
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($sheet2021)
$newSheet.Name = "2022-Q3"

# NOTE: the $sheet2021 handle captured above now resolves to whatever
# sheet sits at that same tab position, which after the insert above is
# the *new* sheet, not "2021-Q4" anymore. Re-resolve it by name so every
# later reference to "the 2021-Q4 sheet" is actually correct.
$sheet2021 = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 2. Populate "2022-Q3" with the fund holdings table
# ---------------------------------------------------------------------
# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet "B2" "100055"
$newSheet.Range("C2").Value = "富国全球科技互联网股票（QDII）"
Set-TextValue $newSheet "D2" "3.95"
Set-TextValue $newSheet "E2" "86.97"
Set-TextValue $newSheet "F2" "3.71"
Set-TextValue $newSheet "G2" "0.1465"
$newSheet.Range("H2").Value = 9

# Row 3
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet "B3" "012920"
$newSheet.Range("C3").Value = "易方达全球成长精选混合（QDII）A（人民币份额）"
Set-TextValue $newSheet "D3" "3.11"
Set-TextValue $newSheet "E3" "68.23"
Set-TextValue $newSheet "F3" "4.00"
Set-TextValue $newSheet "G3" "0.1244"
$newSheet.Range("H3").Value = 2

# Row 4
$newSheet.Range("A4").Value = 2
Set-TextValue $newSheet "B4" "012921"
$newSheet.Range("C4").Value = "易方达全球成长精选混合（QDII）A（美元现汇份额）"
Set-TextValue $newSheet "D4" "3.11"
Set-TextValue $newSheet "E4" "68.23"
Set-TextValue $newSheet "F4" "4.00"
Set-TextValue $newSheet "G4" "0.1244"
$newSheet.Range("H4").Value = 2

# Row 5
$newSheet.Range("A5").Value = 3
Set-TextValue $newSheet "B5" "012922"
$newSheet.Range("C5").Value = "易方达全球成长精选混合（QDII）C（人民币份额）"
Set-TextValue $newSheet "D5" "3.11"
Set-TextValue $newSheet "E5" "68.23"
Set-TextValue $newSheet "F5" "4.00"
Set-TextValue $newSheet "G5" "0.1244"
$newSheet.Range("H5").Value = 2

# Row 6
$newSheet.Range("A6").Value = 4
Set-TextValue $newSheet "B6" "012923"
$newSheet.Range("C6").Value = "易方达全球成长精选混合（QDII）C（美元现汇份额）"
Set-TextValue $newSheet "D6" "3.11"
Set-TextValue $newSheet "E6" "68.23"
Set-TextValue $newSheet "F6" "4.00"
Set-TextValue $newSheet "G6" "0.1244"
$newSheet.Range("H6").Value = 2

# Row 7
$newSheet.Range("A7").Value = 5
Set-TextValue $newSheet "B7" "006555"
$newSheet.Range("C7").Value = "浦银安盛全球智能科技股票（QDII）A"
Set-TextValue $newSheet "D7" "0.25"
Set-TextValue $newSheet "E7" "84.65"
Set-TextValue $newSheet "F7" "6.98"
Set-TextValue $newSheet "G7" "0.0174"
$newSheet.Range("H7").Value = 3

# Row 8
$newSheet.Range("A8").Value = 6
Set-TextValue $newSheet "B8" "014002"
$newSheet.Range("C8").Value = "浦银安盛全球智能科技股票（QDII）C"
Set-TextValue $newSheet "D8" "0.01"
Set-TextValue $newSheet "E8" "84.65"
Set-TextValue $newSheet "F8" "6.98"
Set-TextValue $newSheet "G8" "0.0007"
$newSheet.Range("H8").Value = 3

# Apply the bold/centered/bordered header style (same style already used
# for the "总计" / "2021-Q4" header rows and A-column index cells) to the
# new sheet's header row and index column.
$sheet2021.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$sheet2021.Range("A2").Copy()
$newSheet.Range("A2:A8").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Update the "总计" summary sheet: insert the 2022-Q3 totals into row 2
#    (reusing the existing A2 styled cell) and push the former row-2
#    ("2021-Q4" totals) down to row 3, copying A2's style to the new A3.
# ---------------------------------------------------------------------
$sheetTotal.Range("A2").Copy()
$sheetTotal.Range("A3").PasteSpecial(-4122)

$sheetTotal.Range("A3").Value = 1
$sheetTotal.Range("B3").Value = "2021-Q4"
$sheetTotal.Range("C3").Value = 1
$sheetTotal.Range("D3").Value = 0.09

$sheetTotal.Range("B2").Value = "2022-Q3"
$sheetTotal.Range("C2").Value = 7
$sheetTotal.Range("D2").Value = 0.66
